$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading percent results for the 380 kV case
$dataBF = @(
    @(17.2740012113214, 6.210456982042762, 3.429573696358188, 10.367439548968, 58.69465117415611),
    @(17.24835937570876, 6.094336099605056, 3.441385103338878, 10.38405863763133, 58.51140703389704),
    @(17.23751755836731, 6.023962917165284, 3.449260142005552, 10.39496118729553, 58.40547193870997),
    @(17.23433732991699, 5.995561236524317, 3.452625987723473, 10.39958015578731, 58.36397429323817),
    @(17.23388414955419, 5.990863088341937, 3.453194351241381, 10.40035778165114, 58.35718510085081),
    @(17.23746965055336, 6.02357870773417, 3.449304900349981, 10.39502276668652, 58.40490549582052),
    @(17.26414585971563, 6.170252342280196, 3.433517069753441, 10.3730251553016, 58.63011452670086),
    @(17.35507739569229, 6.463332278339619, 3.407494370419144, 10.33540756791361, 59.12310140208095),
    @(17.44498017893455, 6.679497485950256, 3.391379461507336, 10.31110540338158, 59.51540454765183),
    @(17.49077490609786, 6.777515173490818, 3.384699404859465, 10.3007677707533, 59.70012591172892),
    @(17.50880861860708, 6.814546539023253, 3.38226331422718, 10.29695587723504, 59.77094877830146),
    @(17.50489415074417, 6.806575608018897, 3.382783812509232, 10.29777227410158, 59.75565743404545),
    @(17.49224474945821, 6.780563709441695, 3.384497112692897, 10.30045210734175, 59.70593521884694),
    @(17.48458640183339, 6.76461833502141, 3.385558733194522, 10.30210694988429, 59.67559176853317),
    @(17.44208493367571, 6.673081939490847, 3.391829106863139, 10.31179539324963, 59.50345624164324),
    @(17.41725782194991, 6.616816928571979, 3.39584238944382, 10.31792240659434, 59.39944111234757),
    @(17.40343968805036, 6.584427242369404, 3.398211974221171, 10.32151406723761, 59.34020536981552),
    @(17.39884076222018, 6.573457177824467, 3.39902479681084, 10.32274175951841, 59.32025148493773),
    @(17.41985300312335, 6.622809582802623, 3.39540882992884, 10.31726318607156, 59.41045272742949),
    @(17.49594149893249, 6.788206674344781, 3.383991337831956, 10.29966219040245, 59.72051637314662),
    @(17.54969834477575, 6.895785833737572, 3.377074314416471, 10.28875758978415, 59.92823520880132),
    @(17.52064292510169, 6.83842893706615, 3.380716218158122, 10.29452294610883, 59.81691682760209),
    @(17.41867830255279, 6.620100432135234, 3.395604648155319, 10.31756100422241, 59.40547261734856),
    @(17.32638519142521, 6.383723089846399, 3.414006243166851, 10.34499626285957, 58.98435577580333)
)

$dataIM = @(
    @(42.75315948101491, 10.42460092597664, 15.86322915611633, 11.72837058400546, 17.62854016503063),
    @(42.65336659491722, 10.43262268719585, 15.84653044710395, 11.75000770832539, 17.65987514556383),
    @(42.59598130303159, 10.43789362703449, 15.84037225040109, 11.76472790485008, 17.68216696656818),
    @(42.57358143914718, 10.4401286812294, 15.83889566593751, 11.77108772213461, 17.69201876254439),
    @(42.5699216609609, 10.44050507737537, 15.83871294526435, 11.77216559296533, 17.69370101781995),
    @(42.59567521273932, 10.43792341679539, 15.84034815072141, 11.76481221250075, 17.68229672272646),
    @(42.71794860964929, 10.4272952392535, 15.85662373624067, 11.73553348978019, 17.6387111963775),
    @(42.98822675007963, 10.40918564381834, 15.92085323074491, 11.68948583016186, 17.57744292963109),
    @(43.20489629359475, 10.39753311512838, 15.98745150769677, 11.66256086518133, 17.54716116971183),
    @(43.30727481102847, 10.39258816344311, 16.02188424569289, 11.65180628615569, 17.5365769525681),
    @(43.34657993512393, 10.39076659405804, 16.03550990821626, 11.64794813338434, 17.5330270142215),
    @(43.33809120436383, 10.39115663748647, 16.03254942098203, 11.64876952825326, 17.53377119819415),
    @(43.31049778346802, 10.39243728129028, 16.02299352761651, 11.65148457974066, 17.53627572158513),
    @(43.29366553324999, 10.39322834553187, 16.01721642165659, 11.65317553198794, 17.5378694431982),
    @(43.19828136973947, 10.39786342326205, 15.98528381346666, 11.66329372403586, 17.54791702097937),
    @(43.14073497908485, 10.40079789517526, 15.96674819989594, 11.66988319060026, 17.55489767630933),
    @(43.1079956061597, 10.40251923081239, 15.95647692356771, 11.67381390015908, 17.55921318333104),
    @(43.09697276356898, 10.40310780590628, 15.9530664497516, 11.6751689355304, 17.56072595927922),
    @(43.14682373796516, 10.40048204939107, 15.96868104894794, 11.66916717971464, 17.55412348678635),
    @(43.31858817496047, 10.39205974320742, 16.02578447020912, 11.65068128874765, 17.53552765740841),
    @(43.43396791521264, 10.38685232054155, 16.06652045058902, 11.63984901017604, 17.52604374948168),
    @(43.37210603590728, 10.38960450421787, 16.04446925705928, 11.64551623279457, 17.5308615272058),
    @(43.14406993496133, 10.40062473651131, 15.96780600771604, 11.66949044483038, 17.55447255634995),
    @(42.91189267283753, 10.41379361887622, 15.9000447374741, 11.70072834054299, 17.59142838742796)
)

for ($i = 0; $i -lt $dataBF.Length; $i++) {
    $row = 2 + $i
    $rowValsBF = $dataBF[$i]
    $rowValsIM = $dataIM[$i]
    for ($j = 0; $j -lt $rowValsBF.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowValsBF[$j]
    }
    for ($j = 0; $j -lt $rowValsIM.Length; $j++) {
        $col = 9 + $j
        $ws.Cells.Item($row, $col).Value = $rowValsIM[$j]
    }
}
